# Adding a function to apply TableStyle to Excel module
# Populate header rows on all three sheets and turn SpaceCheck's header
# range into a styled Excel Table ("Tableau1"), matching the
# ImportExportTemplate.xlsx authoring workflow.

$wb = $excel.ActiveWorkbook

# --- SpaceCheck sheet (table with header + TableStyleMedium15) ---
$wsCheck = $wb.Worksheets.Item("SpaceCheck")
$wsCheck.Range("A1").Value = "Service"
$wsCheck.Range("B1").Value = "Local"
$wsCheck.Range("C1").Value = "Remarque"

$tbl = $wsCheck.ListObjects.Add(1, $wsCheck.Range("A1:C14"), $null, 1)
$tbl.Name = "Tableau1"
$tbl.TableStyle = "TableStyleMedium15"

$wsCheck.Columns.Item(1).ColumnWidth = 10
$wsCheck.Columns.Item(2).ColumnWidth = 10
$wsCheck.Columns.Item(3).ColumnWidth = 65

$wsCheck.PageSetup.PaperSize = 9
$wsCheck.PageSetup.Orientation = 1

$wsCheck.Range("C21").Select()

# --- Import sheet ---
$wsImport = $wb.Worksheets.Item("Import")
$wsImport.Range("A1").Value = "Service"
$wsImport.Range("B1").Value = "Local"

$wsImport.PageSetup.PaperSize = 9
$wsImport.PageSetup.Orientation = 1

$wsImport.Range("A1:B1").Select()

# --- Export sheet ---
$wsExport = $wb.Worksheets.Item("Export")
$wsExport.Range("A1").Value = "Service"
$wsExport.Range("B1").Value = "Local"

$wsExport.Range("B10").Select()

# Re-select the SpaceCheck sheet/cell so it remains the active tab.
$wsCheck.Activate()
$wsCheck.Range("C21").Select()
